$d = $word.ActiveDocument

# 1. Replace the long, stitched-together "MANDATORY Fusion Leader Meeting ..."
#    sentence (spread across several runs) with the new single sentence.
#    Word's Find/Replace merges the matched span into one run, which mirrors
#    the diff (six runs collapse into the one new run).
$old = ": there will be a MANDATORY Fusion Leader Meeting tonight where we will discuss Fusion plans & expectations and answer all your questions."
$new = ": We will have a mandatory leader meeting tonight at 7:00 PM in the Chapel while students are in Large Group. We will discuss Fusion rules and expectations, as well as answer all your questions."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 2. Locate the "Attention All Leaders" paragraph (now ending in the
#    replaced sentence) so the trailing _GoBack bookmark can be moved there.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Attention All Leaders*") {
        $target = $p
    }
}

$r = $target.Range

# 3. Move the _GoBack bookmark from the "Close small group..." paragraph to
#    the very end of this paragraph's text (right after the last run, before
#    the paragraph mark). A truly collapsed Range positioned exactly at
#    "end of paragraph text" confuses Bookmarks.Add in this host, so a
#    one-character placeholder is inserted, bookmarked, and then its text is
#    removed again -- leaving a clean, adjacent bookmarkStart/bookmarkEnd
#    pair with no extra characters.
$insertPos = $r.End - 1
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.Text = "X"

$bookmarkRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$cleanup = $d.Range($insertPos, $insertPos + 1)
$cleanup.Text = ""
